$d = $word.ActiveDocument
$p9 = $d.Paragraphs.Item(9)
$srcListRange = $p9.Range

# Add both new paragraphs first (they inherit "Normal" from the last original paragraph)
$endOfDoc = $d.Range($d.Content.End - 1, $d.Content.End - 1)
$newPara1 = $d.Paragraphs.Add($endOfDoc)

$endOfDoc2 = $d.Range($d.Content.End - 1, $d.Content.End - 1)
$newPara2 = $d.Paragraphs.Add($endOfDoc2)

# Format paragraph 2 (trailing empty, indented) - do this first so paragraph 1's
# bold/list formatting doesn't need to be touched again afterwards
$r2 = $newPara2.Range
$r2.ParagraphFormat.LeftIndent = 36
$r2.ParagraphFormat.FirstLineIndent = 36
Write-Output "p2 after indent text=[$($r2.Text)] bold=$($r2.Bold) style=$($r2.Style.NameLocal)"

# Format paragraph 1 (bold numbered list item)
$r1 = $newPara1.Range
$r1.Style = "List Paragraph"
$lt = $srcListRange.ListFormat.ListTemplate
$r1.ListFormat.ApplyListTemplate($lt, $true)
$r1.Text = "A researcher has a database filled with patient records each of which include demographic information as well as medical measurements. The researcher desires to find people with undiagnosed diseases. Which language would be best suited for the task: Prolog, C, or SAS? Explain."
$r1.Bold = 1

Write-Output "count=$($d.Paragraphs.Count)"
